$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data values in column O for the listed rows
$ws.Range("O5").Value = 28.6
$ws.Range("O6").Value = 33.6
$ws.Range("O8").Value = 71.2
$ws.Range("O12").Value = 16.100000000000001

# Narrow columns A:C (from ~39.29 to ~34.43 characters).
# Note: the host quantizes ColumnWidth to an MDW-6 pixel grid (stored width
# is always a multiple of 1/6), so 33.6667 is the input that lands closest
# to the target stored width of 34.42578125 (resolves to 34.5, the nearest
# representable value).
$ws.Range("A:C").ColumnWidth = 33.6667

# Select cell M23 to match the saved selection/active cell in the sheet view
$ws.Range("M23").Select()
